$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 5.75
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("AI4").Value = 26
$ws.Range("AK4").Value = 67
$ws.Range("AM4").Value = 67
$ws.Range("AQ4").Value = 34
$ws.Range("AS4").Value = 301
$ws.Range("AX4").Value = 34
$ws.Range("BA4").Value = 201
$ws.Range("G6").Value = 1.33
$ws.Range("H6").Value = 4.5
$ws.Range("I6").Value = 7.3
$ws.Range("J6").Value = 1.75
$ws.Range("L6").Value = 6.7
$ws.Range("P6").Value = 4.1
$ws.Range("Q6").Value = 1.57
$ws.Range("R6").Value = 2.12
$ws.Range("X6").Value = 5.8
$ws.Range("Y6").Value = 7.2
$ws.Range("Z6").Value = 7.3
$ws.Range("AA6").Value = 9
$ws.Range("AD6").Value = 8
$ws.Range("AE6").Value = 16.5
$ws.Range("AF6").Value = 65
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 19
$ws.Range("AK6").Value = 120
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 55
$ws.Range("AN6").Value = 3.2
$ws.Range("AO6").Value = 5.8
$ws.Range("AQ6").Value = 15.5
$ws.Range("AS6").Value = 200
$ws.Range("AT6").Value = 3.15
$ws.Range("AU6").Value = 8.25
$ws.Range("AV6").Value = 75
$ws.Range("AW6").Value = 8.75
$ws.Range("AX6").Value = 45
$ws.Range("AY6").Value = 45
$ws.Range("AZ6").Value = 300
$ws.Range("BA6").Value = 300
$ws.Range("G7").Value = 2.92
$ws.Range("H7").Value = 2.92
$ws.Range("I7").Value = 2.47
$ws.Range("J7").Value = 3.5
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 3.1
$ws.Range("M7").Value = 1.02
$ws.Range("N7").Value = 8.8
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 1.91
$ws.Range("R7").Value = 1.8
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.45
$ws.Range("U7").Value = 1.62
$ws.Range("V7").Value = 2.05
$ws.Range("W7").Value = 9.25
$ws.Range("X7").Value = 16
$ws.Range("Z7").Value = 40
$ws.Range("AA7").Value = 25
$ws.Range("AB7").Value = 30
$ws.Range("AC7").Value = 9
$ws.Range("AD7").Value = 5.7
$ws.Range("AE7").Value = 11.75
$ws.Range("AF7").Value = 50
$ws.Range("AG7").Value = 350
$ws.Range("AH7").Value = 8.5
$ws.Range("AI7").Value = 13
$ws.Range("AK7").Value = 29
$ws.Range("AL7").Value = 20
$ws.Range("AM7").Value = 26
$ws.Range("AN7").Value = 4.85
$ws.Range("AO7").Value = 16.5
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 80
$ws.Range("AR7").Value = 110
$ws.Range("AS7").Value = 250
$ws.Range("AT7").Value = 2.42
$ws.Range("AU7").Value = 6.5
$ws.Range("AV7").Value = 55
$ws.Range("AW7").Value = 4.4
$ws.Range("AX7").Value = 13.5
$ws.Range("AY7").Value = 20
$ws.Range("AZ7").Value = 60
$ws.Range("BA7").Value = 90
$ws.Range("BB7").Value = 250
$ws.Range("G9").Value = 1.6
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 5.5
$ws.Range("J9").Value = 2.1
$ws.Range("K9").Value = 2.38
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 2.15
$ws.Range("W9").Value = 8
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 12
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 13
$ws.Range("AD9").Value = 7.5
$ws.Range("AG9").Value = 201
$ws.Range("AH9").Value = 17
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 17
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 8
$ws.Range("AQ9").Value = 23
$ws.Range("AW9").Value = 7
$ws.Range("AX9").Value = 26
$ws.Range("AZ9").Value = 81
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.5
$ws.Range("J10").Value = 2.75
$ws.Range("L10").Value = 4
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 1.85
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 1.91
$ws.Range("W10").Value = 7.5
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 9
$ws.Range("Z10").Value = 19
$ws.Range("AA10").Value = 17
$ws.Range("AD10").Value = 6.5
$ws.Range("AI10").Value = 17
$ws.Range("AJ10").Value = 12
$ws.Range("AL10").Value = 29
$ws.Range("AM10").Value = 34
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 11
$ws.Range("AQ10").Value = 41
$ws.Range("AW10").Value = 5.5
$ws.Range("AX10").Value = 19
$ws.Range("AZ10").Value = 67
$ws.Range("BA10").Value = 81
$ws.Range("BC10").Value = 151
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 2.88
$ws.Range("K11").Value = 1.91
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.57
$ws.Range("P11").Value = 2.25
$ws.Range("Q11").Value = 2.7
$ws.Range("R11").Value = 1.44
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 2.2
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53
$ws.Range("W11").Value = 5
$ws.Range("X11").Value = 8
$ws.Range("Y11").Value = 10
$ws.Range("AA11").Value = 21
$ws.Range("AB11").Value = 41
$ws.Range("AC11").Value = 6
$ws.Range("AE11").Value = 21
$ws.Range("AF11").Value = 81
$ws.Range("AH11").Value = 8
$ws.Range("AI11").Value = 19
$ws.Range("AJ11").Value = 15
$ws.Range("AL11").Value = 41
$ws.Range("AM11").Value = 51
$ws.Range("AN11").Value = 3.75
$ws.Range("AP11").Value = 29
$ws.Range("AR11").Value = 81
$ws.Range("AS11").Value = 301
$ws.Range("AT11").Value = 2.2
$ws.Range("AU11").Value = 10
$ws.Range("AV11").Value = 81
$ws.Range("AX11").Value = 26
$ws.Range("AY11").Value = 41
$ws.Range("AZ11").Value = 101
$ws.Range("BA11").Value = 151
$ws.Range("BB11").Value = 351
$ws.Range("G13").Value = 1.36
$ws.Range("H13").Value = 4.33
$ws.Range("I13").Value = 9.5
$ws.Range("K13").Value = 2.4
$ws.Range("L13").Value = 8
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.98
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("W13").Value = 6.5
$ws.Range("Z13").Value = 8.5
$ws.Range("AA13").Value = 12
$ws.Range("AD13").Value = 8.5
$ws.Range("AE13").Value = 21
$ws.Range("AF13").Value = 67
$ws.Range("AK13").Value = 101
$ws.Range("AW13").Value = 9
$ws.Range("BA13").Value = 201
